# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): set as Text to preserve formatted strings
#     (e.g. thousands separated by dots, trailing zeros) exactly as scraped.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.643.79"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.556.74"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.34"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.93"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.557.85"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.86"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.132"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.411"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.151.59"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000197"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.75"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.539.52"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.419.69"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.15"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.70"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "424.78"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.598"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.39"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.695.94"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000117"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.98"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.48"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.95"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.560.42"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.156"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.36"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.66"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.66"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "175.32"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.30"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0824"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.47"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.78"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.89"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.10"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.916"

# --- Volume(1h) column (E): plain text percentage strings, no numeric coercion risk.
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +5.42%  "
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E14").Value = "  -5.74%  "
$ws.Range("E15").Value = "  -4.04%  "
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("E23").Value = "  -3.48%  "
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -7.09%  "
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  -5.82%  "
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("E45").Value = "  -6.95%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  -5.68%  "
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("E51").Value = "  -3.63%  "

# --- Rows 42/43: Mantle and Filecoin swapped position in the ranking,
#     each also getting a refreshed price and 1h volume.
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.01"
$ws.Range("E42").Value = "  -4.26%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.864"
$ws.Range("E43").Value = "  -3.71%  "
